$wb = $excel.ActiveWorkbook

# --- Sheet: weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.70892943143692
$ws.Range("C2").Value = 0.395908632478748
$ws.Range("B3").Value = 0.242382907798554
$ws.Range("C3").Value = 0.219530970391322

# --- Sheet: lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.89257104125165
$ws.Range("C2").Value = 0.410670123606949
$ws.Range("B3").Value = -0.987198951461175
$ws.Range("C3").Value = 0.173663665314817

# --- Sheet: llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.94606176905585
$ws.Range("C2").Value = 0.208222229400172
$ws.Range("B3").Value = 0.536281630291043
$ws.Range("C3").Value = 0.18904560480009

# --- Sheet: gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.50741289770736
$ws.Range("C2").Value = 0.318208850704761
$ws.Range("B3").Value = 0.0248428547262417
$ws.Range("C3").Value = 0.047688533074933

# --- Sheet: exp (no value changes) ---

# --- Sheet: weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.156743645271192
$ws.Range("B2").Value = -0.0662510125415888
$ws.Range("A3").Value = -0.0662510125415888
$ws.Range("B3").Value = 0.0481938469609553

# --- Sheet: lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.168649950423347
$ws.Range("B2").Value = -0.0611261870872031
$ws.Range("A3").Value = -0.0611261870872031
$ws.Range("B3").Value = 0.0301590686505767

# --- Sheet: llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0433564968163777
$ws.Range("B2").Value = -0.00355768253319537
$ws.Range("A3").Value = -0.00355768253319537
$ws.Range("B3").Value = 0.0357382406942316

# --- Sheet: gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.101256872666845
$ws.Range("B2").Value = -0.0101235506567224
$ws.Range("A3").Value = -0.0101235506567224
$ws.Range("B3").Value = 0.00227419618683898

# --- Sheet: exp cov (no value changes) ---

$wb.Save()
